# Insert a new data row at row 160 (shifts existing rows 160-189 down to 161-190)
# and populate it with the new weekly record, per the commit's diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(160).Insert($xlShiftDown)

$ws.Cells.Item(160, 1).Value = 11
$ws.Cells.Item(160, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(160, 3).Value = "Bíobío"
$ws.Cells.Item(160, 4).Value = 45034
$ws.Cells.Item(160, 5).Value = 8
$ws.Cells.Item(160, 6).Value = 100112021
$ws.Cells.Item(160, 7).Value = "Ají"
$ws.Cells.Item(160, 8).Value = "Americana (o)"
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 100
$ws.Cells.Item(160, 11).Value = 11000
$ws.Cells.Item(160, 12).Value = 12000
$ws.Cells.Item(160, 13).Value = 11500
$ws.Cells.Item(160, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(160, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(160, 16).Value = 460
$ws.Cells.Item(160, 17).Value = 25
$ws.Cells.Item(160, 18).Value = "Hortaliza"
